$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

$ws.Cells.Item(2, 1).Value = "와이바이오로직스"
$ws.Cells.Item(2, 2).Value = "2023.11.10~11.16"
$ws.Cells.Item(2, 3).Value = "9,000~11,000"
$ws.Cells.Item(2, 4).Value = "-"
$ws.Cells.Item(2, 5).Value = 13500
$ws.Cells.Item(2, 6).Value = "유안타증권"

$ws.Cells.Item(3, 1).Value = "한선엔지니어링"
$ws.Cells.Item(3, 2).Value = "2023.11.02~11.08"
$ws.Cells.Item(3, 3).Value = "5,200~6,000"
$ws.Cells.Item(3, 4).Value = "-"
$ws.Cells.Item(3, 5).Value = 22100
$ws.Cells.Item(3, 6).Value = "대신증권"

$ws.Cells.Item(4, 1).Value = "에코아이"
$ws.Cells.Item(4, 2).Value = "2023.11.01~11.07"
$ws.Cells.Item(4, 3).Value = "28,500~34,700"
$ws.Cells.Item(4, 4).Value = "-"
$ws.Cells.Item(4, 5).Value = 59251
$ws.Cells.Item(4, 6).Value = "KB증권"

$ws.Cells.Item(5, 1).Value = "동인기연(유가)"
$ws.Cells.Item(5, 2).Value = "2023.11.01~11.07"
$ws.Cells.Item(5, 3).Value = "33,000~37,000"
$ws.Cells.Item(5, 4).Value = "-"
$ws.Cells.Item(5, 5).Value = 60654
$ws.Cells.Item(5, 6).Value = "NH투자증권"

$ws.Cells.Item(6, 1).Value = "스톰테크"
$ws.Cells.Item(6, 2).Value = "2023.10.31~11.06"
$ws.Cells.Item(6, 3).Value = "8,000~9,500"
$ws.Cells.Item(6, 4).Value = "-"
$ws.Cells.Item(6, 5).Value = 26800
$ws.Cells.Item(6, 6).Value = "하이투자증권"

$ws.Cells.Item(7, 1).Value = "블루엠텍"
$ws.Cells.Item(7, 2).Value = "2023.10.31~11.06"
$ws.Cells.Item(7, 3).Value = "15,000~19,000"
$ws.Cells.Item(7, 4).Value = "-"
$ws.Cells.Item(7, 5).Value = 21000
$ws.Cells.Item(7, 6).Value = "하나증권,키움증권"

$ws.Cells.Item(8, 1).Value = "에코프로머티리얼즈"
$ws.Cells.Item(8, 2).Value = "2023.10.30~11.03"
$ws.Cells.Item(8, 3).Value = "36,200~44,000"
$ws.Cells.Item(8, 4).Value = "-"
$ws.Cells.Item(8, 5).Value = 524031
$ws.Cells.Item(8, 6).Value = "미래에셋증권,NH투자증권,하이투자증권"

$ws.Cells.Item(9, 1).Value = "캡스톤파트너스"
$ws.Cells.Item(9, 2).Value = "2023.10.26~11.01"
$ws.Cells.Item(9, 3).Value = "3,200~3,600"
$ws.Cells.Item(9, 4).Value = "-"
$ws.Cells.Item(9, 5).Value = 5107
$ws.Cells.Item(9, 6).Value = "NH투자증권"

$ws.Cells.Item(10, 1).Value = "에이텀"
$ws.Cells.Item(10, 2).Value = "2023.10.26~11.01"
$ws.Cells.Item(10, 3).Value = "23,000~30,000"
$ws.Cells.Item(10, 4).Value = "-"
$ws.Cells.Item(10, 5).Value = 14950
$ws.Cells.Item(10, 6).Value = "하나증권"

$ws.Cells.Item(11, 1).Value = "한국스팩13호"
$ws.Cells.Item(11, 2).Value = "2023.10.25~10.26"
$ws.Cells.Item(11, 3).Value = "2,000~2,000"
$ws.Cells.Item(11, 4).Value = "-"
$ws.Cells.Item(11, 5).Value = 8000
$ws.Cells.Item(11, 6).Value = "한국투자증권"

$ws.Cells.Item(12, 1).Value = "그린리소스"
$ws.Cells.Item(12, 2).Value = "2023.10.25~10.31"
$ws.Cells.Item(12, 3).Value = "11,000~14,000"
$ws.Cells.Item(12, 4).Value = "-"
$ws.Cells.Item(12, 5).Value = 18040
$ws.Cells.Item(12, 6).Value = "NH투자증권"

$ws.Cells.Item(13, 1).Value = "에이직랜드"
$ws.Cells.Item(13, 2).Value = "2023.10.23~10.27"
$ws.Cells.Item(13, 3).Value = "19,100~21,400"
$ws.Cells.Item(13, 4).Value = "-"
$ws.Cells.Item(13, 5).Value = 50353
$ws.Cells.Item(13, 6).Value = "삼성증권"

$ws.Cells.Item(14, 1).Value = "에스와이스틸텍"
$ws.Cells.Item(14, 2).Value = "2023.10.23~10.27"
$ws.Cells.Item(14, 3).Value = "1,200~1,500"
$ws.Cells.Item(14, 4).Value = "-"
$ws.Cells.Item(14, 5).Value = 8400
$ws.Cells.Item(14, 6).Value = "KB증권"

$ws.Cells.Item(15, 1).Value = "컨텍"
$ws.Cells.Item(15, 2).Value = "2023.10.20~10.26"
$ws.Cells.Item(15, 3).Value = "20,300~22,500"
$ws.Cells.Item(15, 4).Value = "-"
$ws.Cells.Item(15, 5).Value = 41818
$ws.Cells.Item(15, 6).Value = "대신증권"

$ws.Cells.Item(16, 1).Value = "큐로셀"
$ws.Cells.Item(16, 2).Value = "2023.10.20~10.26"
$ws.Cells.Item(16, 3).Value = "29,800~33,500"
$ws.Cells.Item(16, 4).Value = "-"
$ws.Cells.Item(16, 5).Value = 47680
$ws.Cells.Item(16, 6).Value = "미래에셋증권,삼성증권"

$ws.Cells.Item(17, 1).Value = "메가터치"
$ws.Cells.Item(17, 2).Value = "2023.10.20~10.26"
$ws.Cells.Item(17, 3).Value = "3,500~4,000"
$ws.Cells.Item(17, 4).Value = "-"
$ws.Cells.Item(17, 5).Value = 18200
$ws.Cells.Item(17, 6).Value = "NH투자증권"

$ws.Cells.Item(18, 1).Value = "비아이매트릭스"
$ws.Cells.Item(18, 2).Value = "2023.10.19~10.25"
$ws.Cells.Item(18, 3).Value = "9,100~11,000"
$ws.Cells.Item(18, 4).Value = "-"
$ws.Cells.Item(18, 5).Value = 10920
$ws.Cells.Item(18, 6).Value = "IBK투자증권"

$ws.Cells.Item(19, 1).Value = "KB스팩27호"
$ws.Cells.Item(19, 2).Value = "2023.10.19~10.20"
$ws.Cells.Item(19, 3).Value = "-"
$ws.Cells.Item(19, 4).Value = "-"
$ws.Cells.Item(19, 5).Value = 25000
$ws.Cells.Item(19, 6).Value = "KB증권"

$ws.Cells.Item(20, 1).Value = "유투바이오"
$ws.Cells.Item(20, 2).Value = "2023.10.18~10.19"
$ws.Cells.Item(20, 3).Value = "3,300~3,900"
$ws.Cells.Item(20, 4).Value = "-"
$ws.Cells.Item(20, 5).Value = 3724
$ws.Cells.Item(20, 6).Value = "신한투자증권"

$ws.Cells.Item(21, 1).Value = "쏘닉스"
$ws.Cells.Item(21, 2).Value = "2023.10.17~10.23"
$ws.Cells.Item(21, 3).Value = "5,000~7,000"
$ws.Cells.Item(21, 4).Value = "-"
$ws.Cells.Item(21, 5).Value = 18000
$ws.Cells.Item(21, 6).Value = "KB증권"

